$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '69.375.10'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +2.23%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.388.30'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +1.77%  '
$ws.Range("E4").Value = '  +0.03%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '586.95'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +0.72%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '179.90'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +1.62%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -0.07%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.595'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +0.83%  '
$ws.Range("E9").Value = '  +5.69%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.590'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +1.36%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '48.47'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +2.07%  '
$ws.Range("E12").Value = '  +3.01%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '677.49'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -3.24%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '8.61'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +1.97%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '3.930.34'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +1.60%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '69.395.87'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +2.19%  '
$ws.Range("E17").Value = '  +1.76%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '3.388.63'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +1.55%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '17.62'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +0.72%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '11.27'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +1.83%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '0.903'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +0.78%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '5.42'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +0.38%  '
$ws.Range("E23").Value = '  +0.52%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '103.30'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +3.71%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '3.93'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("E26").Value = '  +1.05%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '9.65'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +0.29%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '34.06'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +2.76%  '
$ws.Range("E29").Value = '  +1.29%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '7.04'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -0.69%  '
$ws.Range("E31").Value = '  +0.89%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '556.06'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -2.04%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '3.61'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +6.25%  '
$ws.Range("E34").Value = '  +0.53%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '58.29'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +1.72%  '
$ws.Range("E36").Value = '  +0.06%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '3.679.75'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -0.21%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.138'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +3.93%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '35.31'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +2.06%  '
$ws.Range("E40").Value = '  +3.05%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '2.71'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +1.65%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0696'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +2.93%  '
$ws.Range("E43").Value = '  +0.41%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.0423'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +4.03%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '3.30'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -1.02%  '
$ws.Range("E46").Value = '  -0.37%  '
$ws.Range("E47").Value = '  +0.73%  '
$ws.Range("E48").Value = '  +5.69%  '
$ws.Range("E49").Value = '  -0.03%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '132.73'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +1.46%  '
$ws.Range("E51").Value = '  +3.61%  '
